$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Apply the "empty underlined cell" formatting first (row 28 / C28) so that
# the resulting style entry (fontId=2: underline, theme color 1) lands
# before any style created later by Hyperlinks.Add.
$ws.Range("C28").Font.Underline = 2

# --- Row 24: themeforest wallet inspiration link ---
$ws.Range("B24").Value = "https://themeforest.net/search/wallet?srsltid=AfmBOoo8oQOEGuS2iThy6Hil9fR1ffhqVDq8GYZmM_M84hnPVaXUbgiR"
$ws.Hyperlinks.Add($ws.Range("B24"), "https://themeforest.net/search/wallet?srsltid=AfmBOoo8oQOEGuS2iThy6Hil9fR1ffhqVDq8GYZmM_M84hnPVaXUbgiR")
$ws.Range("B24").Style = "Hipervínculo"
$ws.Range("C24").Value = "Inspiración para Wallet"

# --- Row 25: dribbble wallet inspiration link ---
$ws.Range("B25").Value = "https://dribbble.com/tags/my-wallet"
$ws.Hyperlinks.Add($ws.Range("B25"), "https://dribbble.com/tags/my-wallet")
$ws.Range("B25").Style = "Hipervínculo"
$ws.Range("C25").Value = "Inspiración para Wallet"

# Match the saved selection state (active cell C28) from the source file.
$ws.Range("C28").Select()
